$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "E" (working) column bucket-list updates ---

# Row 2 (pcb): E2 was "partially" -> now "yes", shown with a 0% number format
# and using the plain (non-"header style") font, matching the new xf record.
$ws.Range("E2").Style = "Normal"
$ws.Range("E2").NumberFormat = "0%"
$ws.Range("E2").Value = "yes"

# Row 3 (kernel): E3 was "partially" -> now a literal 0.5 (displayed as 50%),
# keeping its previous (bold-ish) font style but adding the 0% number format.
$ws.Range("E3").NumberFormat = "0%"
$ws.Range("E3").Value = 0.5

# Row 4 (kernel::lock): E4 was "no" -> now "yes"
$ws.Range("E4").Value = "yes"

# Column E (the new percentages) is now a bit narrower than before
$ws.Columns("E").ColumnWidth = 11.29

# Selection moved from E5 to B10
$ws.Range("B10").Select()
